$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = "e2cbc8ae-1a14-4e4b-b0f7-637666520407.md"
$ws1.Range("B2").Value = "e2e\e2cbc8ae-1a14-4e4b-b0f7-637666520407.md"
$ws1.Range("A3").Value = "7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.md"
$ws1.Range("B3").Value = "e2e\7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.md"
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-18 08:50:35"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = "e2cbc8ae-1a14-4e4b-b0f7-637666520407.md"
$ws2.Range("G2").Value = "e2cbc8ae-1a14-4e4b-b0f7-637666520407.c654b263278e6a1a91fbe95efde0c348a508b69f.zh-cn.xlf"
$ws2.Range("I2").Value = "e2cbc8ae-1a14-4e4b-b0f7-637666520407.md"
$ws2.Range("J2").Value = "e2cbc8ae-1a14-4e4b-b0f7-637666520407.c654b263278e6a1a91fbe95efde0c348a508b69f.zh-cn.xlf"
$ws2.Range("A3").Value = "7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("G3").Value = "7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.e53680d43fa72cb1a8d826856d789cacaacbcb81.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-18 08:50:30"
$ws2.Range("I3").Value = "7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.md"
$ws2.Range("J3").Value = "7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.e53680d43fa72cb1a8d826856d789cacaacbcb81.zh-cn.xlf"
$ws2.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/27f263874f3dbecc576224e1f0447730252cbee0/e2e/7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8bfc4d4734404b7be21769e18f9ea62079e3b35f/e2e/7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.md."
$ws2.Columns.Item(16).ColumnWidth = 40

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = "e2cbc8ae-1a14-4e4b-b0f7-637666520407.md"
$ws3.Range("G2").Value = "e2cbc8ae-1a14-4e4b-b0f7-637666520407.c654b263278e6a1a91fbe95efde0c348a508b69f.de-de.xlf"
$ws3.Range("I2").Value = "e2cbc8ae-1a14-4e4b-b0f7-637666520407.md"
$ws3.Range("J2").Value = "e2cbc8ae-1a14-4e4b-b0f7-637666520407.c654b263278e6a1a91fbe95efde0c348a508b69f.de-de.xlf"
$ws3.Range("A3").Value = "7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("G3").Value = "7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.e53680d43fa72cb1a8d826856d789cacaacbcb81.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-18 08:50:35"
$ws3.Range("I3").Value = "7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.md"
$ws3.Range("J3").Value = "7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.e53680d43fa72cb1a8d826856d789cacaacbcb81.de-de.xlf"
$ws3.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/27f263874f3dbecc576224e1f0447730252cbee0/e2e/7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8bfc4d4734404b7be21769e18f9ea62079e3b35f/e2e/7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.md."
$ws3.Columns.Item(16).ColumnWidth = 40
